$d = $word.ActiveDocument

$d.Content.Find.Execute("50×90=", $true, $false, $false, $false, $false, $true, 1, $false, "32×27=", 2)
$d.Content.Find.Execute("47×18=", $true, $false, $false, $false, $false, $true, 1, $false, "52×82=", 2)
$d.Content.Find.Execute("84×34=", $true, $false, $false, $false, $false, $true, 1, $false, "13×97=", 2)
$d.Content.Find.Execute("76×80=", $true, $false, $false, $false, $false, $true, 1, $false, "57×19=", 2)
$d.Content.Find.Execute("84×68=", $true, $false, $false, $false, $false, $true, 1, $false, "42×42=", 2)
$d.Content.Find.Execute("60×39=", $true, $false, $false, $false, $false, $true, 1, $false, "74×53=", 2)
$d.Content.Find.Execute("39×76=", $true, $false, $false, $false, $false, $true, 1, $false, "52×72=", 2)
$d.Content.Find.Execute("60×87=", $true, $false, $false, $false, $false, $true, 1, $false, "73×47=", 2)
$d.Content.Find.Execute("53×79=", $true, $false, $false, $false, $false, $true, 1, $false, "67×21=", 2)
$d.Content.Find.Execute("24×58=", $true, $false, $false, $false, $false, $true, 1, $false, "39×11=", 2)
$d.Content.Find.Execute("19×50=", $true, $false, $false, $false, $false, $true, 1, $false, "17×68=", 2)
$d.Content.Find.Execute("80×64=", $true, $false, $false, $false, $false, $true, 1, $false, "40×94=", 2)
$d.Content.Find.Execute("17×13=", $true, $false, $false, $false, $false, $true, 1, $false, "32×29=", 2)
$d.Content.Find.Execute("51×94=", $true, $false, $false, $false, $false, $true, 1, $false, "39×64=", 2)
$d.Content.Find.Execute("74×47=", $true, $false, $false, $false, $false, $true, 1, $false, "99×50=", 2)
$d.Content.Find.Execute("12×65=", $true, $false, $false, $false, $false, $true, 1, $false, "34×93=", 2)
$d.Content.Find.Execute("57×74=", $true, $false, $false, $false, $false, $true, 1, $false, "16×70=", 2)
$d.Content.Find.Execute("51×71=", $true, $false, $false, $false, $false, $true, 1, $false, "22×44=", 2)
$d.Content.Find.Execute("87×46=", $true, $false, $false, $false, $false, $true, 1, $false, "36×51=", 2)
$d.Content.Find.Execute("21×98=", $true, $false, $false, $false, $false, $true, 1, $false, "37×42=", 2)
$d.Content.Find.Execute("82×26=", $true, $false, $false, $false, $false, $true, 1, $false, "47×50=", 2)
$d.Content.Find.Execute("53×74=", $true, $false, $false, $false, $false, $true, 1, $false, "57×13=", 2)
$d.Content.Find.Execute("62×41=", $true, $false, $false, $false, $false, $true, 1, $false, "58×98=", 2)
$d.Content.Find.Execute("76×45=", $true, $false, $false, $false, $false, $true, 1, $false, "53×92=", 2)
$d.Content.Find.Execute("94×23=", $true, $false, $false, $false, $false, $true, 1, $false, "30×30=", 2)
